# Update gh-pages to output generated at 456a3b4
#
# This updates the "想去人数" (F column) attendance counters on the
# 展览 (Worksheets.Item(1)) and 全部类型 (Worksheets.Item(4)) sheets,
# and marks the cancelled event in C3 on both of those sheets.

$wb = $excel.ActiveWorkbook

$cancelledName = "苏州·第一届寒假动漫展宅舞比赛-CF01（取消）"

# ---- Sheet "展览" (index 1) ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 181
$ws1.Range("C3").Value = $cancelledName
$ws1.Range("F3").Value = 286
$ws1.Range("F4").Value = 138
$ws1.Range("F5").Value = 1292
$ws1.Range("F6").Value = 18052
$ws1.Range("F7").Value = 358
$ws1.Range("F8").Value = 257
$ws1.Range("F10").Value = 6806
$ws1.Range("F12").Value = 157
$ws1.Range("F13").Value = 12
$ws1.Range("F16").Value = 18
$ws1.Range("F17").Value = 151
$ws1.Range("F19").Value = 211
$ws1.Range("F23").Value = 30
$ws1.Range("F26").Value = 983
$ws1.Range("F27").Value = 118
$ws1.Range("F29").Value = 533
$ws1.Range("F30").Value = 28
$ws1.Range("F32").Value = 69
$ws1.Range("F33").Value = 12029
$ws1.Range("F34").Value = 1277
$ws1.Range("F36").Value = 204
$ws1.Range("F37").Value = 273
$ws1.Range("F38").Value = 3911

# ---- Sheet "全部类型" (index 4) ----
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 181
$ws4.Range("C3").Value = $cancelledName
$ws4.Range("F3").Value = 286
$ws4.Range("F4").Value = 138
$ws4.Range("F5").Value = 1292
$ws4.Range("F6").Value = 18052
$ws4.Range("F7").Value = 358
$ws4.Range("F8").Value = 257
$ws4.Range("F10").Value = 6806
$ws4.Range("F12").Value = 157
$ws4.Range("F13").Value = 12
$ws4.Range("F16").Value = 18
$ws4.Range("F17").Value = 151
$ws4.Range("F19").Value = 211
$ws4.Range("F23").Value = 30
$ws4.Range("F26").Value = 983
$ws4.Range("F27").Value = 118
$ws4.Range("F29").Value = 533
$ws4.Range("F32").Value = 28
$ws4.Range("F34").Value = 69
$ws4.Range("F35").Value = 12029
$ws4.Range("F36").Value = 1277
$ws4.Range("F38").Value = 204
$ws4.Range("F39").Value = 273
$ws4.Range("F40").Value = 3911
